$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in status for the 2 newly completed functions (rows 5 & 6, column F)
$ws.Range("F5").Value = "DOne"
$ws.Range("F6").Value = "Done"

# Update the active selection to F7
$ws.Range("F7").Select()
